$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BECbIC")

$ws.Range("B2").Value = 79323799.85092574
$ws.Range("C2").Value = 11012645.44913459
$ws.Range("D2").Value = 296521345.48829
$ws.Range("F2").Value = 1118333047.978115
$ws.Range("G2").Value = 308145992.1600693
$ws.Range("H2").Value = 16565817.92246986
$ws.Range("I2").Value = 144007686.3165732
$ws.Range("J2").Value = 253362597.197578
$ws.Range("L2").Value = 1135515931.960934
$ws.Range("M2").Value = 14918901.52884257
$ws.Range("N2").Value = 93685147.16244857
$ws.Range("O2").Value = 24811795.77544174
$ws.Range("P2").Value = 65655069.09551322
$ws.Range("Q2").Value = 54133407.45706175
$ws.Range("S2").Value = 417063638.7192398
$ws.Range("T2").Value = 27502912.87584348
$ws.Range("U2").Value = 36685041.56455737
$ws.Range("V2").Value = 230815628.0157304
$ws.Range("W2").Value = 19449398.85564813
$ws.Range("X2").Value = 155242340.4850947
$ws.Range("Y2").Value = 90237258.89501102
$ws.Range("Z2").Value = 755781306.4975775
$ws.Range("AA2").Value = 375769684.3117697
$ws.Range("AB2").Value = 611180480.6212597
$ws.Range("AC2").Value = 3609204614.539937
$ws.Range("AD2").Value = 1852309068.784548
$ws.Range("AF2").Value = 1716663429.518345
$ws.Range("AG2").Value = 474742432.5695535
$ws.Range("AK2").Value = 636276989.7630531
$ws.Range("AL2").Value = 3676534034.566645
$ws.Range("AN2").Value = 2204024000
$ws.Range("AO2").Value = 15540447000
$ws.Range("AP2").Value = 369912532.1078408
